# Adds the report rows (reservation data) produced by the "relatorios prontos"
# export into the active worksheet, below the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(9,  "2021-10-22T12:02:18-03:00", 4,  8,  3.75),
    @(10, "2021-10-22T12:02:33-03:00", 5,  6,  3.6),
    @(12, "2021-10-22T14:17:36-03:00", 5,  3,  3.6),
    @(13, "2021-10-22T14:40:58-03:00", 4,  8,  3.75),
    @(15, "2021-10-22T16:26:04-03:00", 11, 14, 2.85),
    @(17, "2021-10-25T14:25:27-03:00", 4,  8,  25.35),
    @(8,  "2021-10-26T10:45:31-03:00", 1,  9,  2.55),
    @(7,  "2021-10-26T16:54:24-03:00", 1,  9,  8.6),
    @(21, "2021-10-29T13:09:11-03:00", 4,  3,  0.15),
    @(5,  "2021-10-29T13:11:13-03:00", 4,  8,  37.15),
    @(6,  "2021-10-29T13:11:18-03:00", 5,  6,  36.8)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}
